{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst newValues = [\n  [\n    \"52-31=\",\n    \"69-44=\",\n    \"17+26=\",\n    \"77-71=\",\n    \"46-11=\"\n  ],\n  [\n    \"52+2=\",\n    \"89-47=\",\n    \"89-74=\",\n    \"77-4=\",\n    \"91+2=\"\n  ],\n  [\n    \"56+14=\",\n    \"35-8=\",\n    \"25-1=\",\n    \"65-65=\",\n    \"86-41=\"\n  ],\n  [\n    \"73+11=\",\n    \"58-29=\",\n    \"32-24=\",\n    \"78-9=\",\n    \"42-6=\"\n  ],\n  [\n    \"29+25=\",\n    \"30+66=\",\n    \"78-69=\",\n    \"14+38=\",\n    \"80-13=\"\n  ],\n  [\n    \"38+31=\",\n    \"10-0=\",\n    \"95-89=\",\n    \"39-33=\",\n    \"32+25=\"\n  ],\n  [\n    \"71-15=\",\n    \"0+75=\",\n    \"10+50=\",\n    \"96-78=\",\n    \"41+46=\"\n  ],\n  [\n    \"2+90=\",\n    \"25+20=\",\n    \"46+42=\",\n    \"62-40=\",\n    \"23+11=\"\n  ],\n  [\n    \"60-36=\",\n    \"39-0=\",\n    \"10+22=\",\n    \"13+66=\",\n    \"51+35=\"\n  ],\n  [\n    \"15-1=\",\n    \"70+29=\",\n    \"6+8=\",\n    \"12-6=\",\n    \"70-12=\"\n  ],\n  [\n    \"95-68=\",\n    \"28+62=\",\n    \"82-62=\",\n    \"20+58=\",\n    \"32-30=\"\n  ],\n  [\n    \"97-1=\",\n    \"62-41=\",\n    \"21+27=\",\n    \"38-12=\",\n    \"93-82=\"\n  ],\n  [\n    \"7+54=\",\n    \"0+46=\",\n    \"44+35=\",\n    \"10+9=\",\n    \"81-80=\"\n  ],\n  [\n    \"3-2=\",\n    \"86+3=\",\n    \"75+1=\",\n    \"32-14=\",\n    \"50-33=\"\n  ],\n  [\n    \"19+51=\",\n    \"99-1=\",\n    \"11+42=\",\n    \"62+2=\",\n    \"67-35=\"\n  ],\n  [\n    \"38+45=\",\n    \"38+40=\",\n    \"72-57=\",\n    \"52-25=\",\n    \"37+55=\"\n  ],\n  [\n    \"50-9=\",\n    \"62+21=\",\n    \"67-53=\",\n    \"75+11=\",\n    \"16+33=\"\n  ],\n  [\n    \"82-57=\",\n    \"2+6=\",\n    \"54-23=\",\n    \"6+15=\",\n    \"21+31=\"\n  ],\n  [\n    \"75+19=\",\n    \"85+10=\",\n    \"26-10=\",\n    \"29+56=\",\n    \"88-67=\"\n  ],\n  [\n    \"76+16=\",\n    \"20+66=\",\n    \"81-32=\",\n    \"88-84=\",\n    \"13+0=\"\n  ]\n];\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$values = @(\n  \"52-31=\",\n  \"69-44=\",\n  \"17+26=\",\n  \"77-71=\",\n  \"46-11=\",\n  \"52+2=\",\n  \"89-47=\",\n  \"89-74=\",\n  \"77-4=\",\n  \"91+2=\",\n  \"56+14=\",\n  \"35-8=\",\n  \"25-1=\",\n  \"65-65=\",\n  \"86-41=\",\n  \"73+11=\",\n  \"58-29=\",\n  \"32-24=\",\n  \"78-9=\",\n  \"42-6=\",\n  \"29+25=\",\n  \"30+66=\",\n  \"78-69=\",\n  \"14+38=\",\n  \"80-13=\",\n  \"38+31=\",\n  \"10-0=\",\n  \"95-89=\",\n  \"39-33=\",\n  \"32+25=\",\n  \"71-15=\",\n  \"0+75=\",\n  \"10+50=\",\n  \"96-78=\",\n  \"41+46=\",\n  \"2+90=\",\n  \"25+20=\",\n  \"46+42=\",\n  \"62-40=\",\n  \"23+11=\",\n  \"60-36=\",\n  \"39-0=\",\n  \"10+22=\",\n  \"13+66=\",\n  \"51+35=\",\n  \"15-1=\",\n  \"70+29=\",\n  \"6+8=\",\n  \"12-6=\",\n  \"70-12=\",\n  \"95-68=\",\n  \"28+62=\",\n  \"82-62=\",\n  \"20+58=\",\n  \"32-30=\",\n  \"97-1=\",\n  \"62-41=\",\n  \"21+27=\",\n  \"38-12=\",\n  \"93-82=\",\n  \"7+54=\",\n  \"0+46=\",\n  \"44+35=\",\n  \"10+9=\",\n  \"81-80=\",\n  \"3-2=\",\n  \"86+3=\",\n  \"75+1=\",\n  \"32-14=\",\n  \"50-33=\",\n  \"19+51=\",\n  \"99-1=\",\n  \"11+42=\",\n  \"62+2=\",\n  \"67-35=\",\n  \"38+45=\",\n  \"38+40=\",\n  \"72-57=\",\n  \"52-25=\",\n  \"37+55=\",\n  \"50-9=\",\n  \"62+21=\",\n  \"67-53=\",\n  \"75+11=\",\n  \"16+33=\",\n  \"82-57=\",\n  \"2+6=\",\n  \"54-23=\",\n  \"6+15=\",\n  \"21+31=\",\n  \"75+19=\",\n  \"85+10=\",\n  \"26-10=\",\n  \"29+56=\",\n  \"88-67=\",\n  \"76+16=\",\n  \"20+66=\",\n  \"81-32=\",\n  \"88-84=\",\n  \"13+0=\"\n)\n$cols = $tbl.Columns.Count\nfor ($i = 0; $i -lt $values.Length; $i++) {\n  $row = [int]([math]::Floor($i / $cols)) + 1\n  $col = ($i % $cols) + 1\n  $tbl.Cell($row, $col).Range.Text = $values[$i]\n}\n"}
